# Automatische test-sync: 2025-08-13 22:26:50
# Append a new log entry (row 26) to the "Logs" sheet, mirroring the previous
# row's content but with a fresh timestamp, then refresh the conditional
# formatting ranges and the "Dashboard" sheet's summary count so they include
# the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# The last existing data row is row 25; duplicate it into the new row 26.
$lastRow = 25
$newRow  = 26

$ws.Range("A" + $lastRow + ":J" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":J" + $newRow).PasteSpecial()

# The new entry was logged at a later timestamp than the one it was copied from.
$ws.Range("F" + $newRow).Value = "2025-08-13 22:26:49"

# Expand the conditional formatting ranges that previously stopped at row 25
# so they also cover the freshly added row 26.
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $oldRange = $ws.Range($col + "2:" + $col + $lastRow)
    $fcs = $oldRange.FormatConditions
    $count = $fcs.Count()
    for ($i = 1; $i -le $count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($ws.Range($col + "2:" + $col + $newRow))
    }
}

# Update the Dashboard summary count for this category to reflect the new row.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 25
